# Replace the numeric month values in column C ("Mes") of the
# "Principales Aeropuertos de Carga" table with their Spanish
# three-letter abbreviations (e.g. 8 -> "Ago.").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$monthNames = @{
    1  = "Ene."
    2  = "Feb."
    3  = "Mar."
    4  = "Abr."
    5  = "May."
    6  = "Jun."
    7  = "Jul."
    8  = "Ago."
    9  = "Sep."
    10 = "Oct."
    11 = "Nov."
    12 = "Dic."
}

for ($row = 6; $row -le 85; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $monthNumber = [int]$cell.Value()
    $cell.Value = $monthNames[$monthNumber]
}
